$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of daily data to the "Serie" table (date 02-11-2021), keeping
# the new date label as plain text in A3 -- exactly like the existing
# "07-10-2021" entry in A2 -- instead of letting Excel auto-convert the
# dd-mm-yyyy-looking text into a date serial number.
#
# A helper cell (far outside the used range) is used to build the text value
# via a formula, whose result is a text string and is therefore never subject
# to Excel's "looks like a date" auto-conversion. That value is then copied
# and pasted as a value into A3, which preserves it as plain text without
# touching any cell number formats/styles.
$helper = $ws.Range("ZZ1000")
$helper.Formula = '="02-11-2021"'
$helper.Copy()
$ws.Range("A3").PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = $false

$ws.Range("B3").Value = 5.8
$ws.Range("C3").Value = 6.2
$ws.Range("D3").Value = 2.1
$ws.Range("E3").Value = 2.5
